$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 5609.75
$ws.Range("I20").Value = 2999.5
$ws.Range("J20").Value = 8220
$ws.Range("K20").Value = 2999.5
$ws.Range("L20").Value = 8220
$ws.Range("M20").Value = -2769.5
$ws.Range("N20").Value = -8680

$ws.Range("H35").Value = 5609.75
$ws.Range("I35").Value = 2999.5
$ws.Range("J35").Value = 8220
$ws.Range("K35").Value = 2999.5
$ws.Range("L35").Value = 8220
$ws.Range("M35").Value = -2620.5
$ws.Range("N35").Value = -8978

$ws.Range("H55").Value = 374.0909
$ws.Range("I55").Value = 485.85715
$ws.Range("J55").Value = 178.5
$ws.Range("K55").Value = 485.85715
$ws.Range("L55").Value = 178.5
$ws.Range("M55").Value = -271.85715
$ws.Range("N55").Value = -606.5

$ws.Range("H106").Value = 11908674
$ws.Range("I106").Value = 18523270
$ws.Range("J106").Value = 2400.6
$ws.Range("K106").Value = 18523270
$ws.Range("L106").Value = 2400.6
$ws.Range("M106").Value = -18522639

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6663.562
$ws.Range("I32").Value = 5699.3135
$ws.Range("J32").Value = 20002.334
$ws.Range("K32").Value = 5699.3135
$ws.Range("L32").Value = 20002.334
$ws.Range("M32").Value = -5412.3135

$ws.Range("H61").Value = 5229.1875
$ws.Range("I61").Value = 8168.294
$ws.Range("J61").Value = 3617.4194
$ws.Range("K61").Value = 8168.294
$ws.Range("L61").Value = 3617.4194
$ws.Range("M61").Value = -7956.294
$ws.Range("N61").Value = -4041.4194

$ws.Range("H132").Value = 2914.1746
$ws.Range("I132").Value = 2298.2122
$ws.Range("J132").Value = 3591.7334
$ws.Range("K132").Value = 6894.6366
$ws.Range("L132").Value = 10775.2002
$ws.Range("M132").Value = -4364.6366

$ws.Range("H135").Value = 55029.832
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 55029.832
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 55029.832
$ws.Range("N135").Value = -65169.832

$ws.Range("H136").Value = 5229.1875
$ws.Range("I136").Value = 8168.294
$ws.Range("J136").Value = 3617.4194
$ws.Range("K136").Value = 24504.882
$ws.Range("L136").Value = 10852.2582
$ws.Range("M136").Value = -21954.882
$ws.Range("N136").Value = -15952.2582

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 4357.75
$ws.Range("I23").Value = 1506
$ws.Range("J23").Value = 7209.5
$ws.Range("K23").Value = 1506
$ws.Range("L23").Value = 7209.5
$ws.Range("M23").Value = -1223
$ws.Range("N23").Value = -7775.5

$ws.Range("H105").Value = 10903.16
$ws.Range("I105").Value = 17800.385
$ws.Range("J105").Value = 3431.1667
$ws.Range("K105").Value = 17800.385
$ws.Range("L105").Value = 3431.1667
$ws.Range("M105").Value = -16053.385
$ws.Range("N105").Value = -6925.1667

$ws.Range("H107").Value = 1208.909
$ws.Range("I107").Value = 1339.6
$ws.Range("J107").Value = 1100
$ws.Range("K107").Value = 1339.6
$ws.Range("L107").Value = 1100
$ws.Range("M107").Value = 580.4000000000001
$ws.Range("N107").Value = -4940

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 4980.8887
$ws.Range("I94").Value = 5582.5557
$ws.Range("J94").Value = 4680.0557
$ws.Range("K94").Value = 5582.5557
$ws.Range("L94").Value = 4680.0557
$ws.Range("M94").Value = -5131.5557
$ws.Range("N94").Value = -5582.0557

$ws.Range("H134").Value = 3164.861
$ws.Range("I134").Value = 3245.25
$ws.Range("J134").Value = 2883.5
$ws.Range("K134").Value = 9735.75
$ws.Range("L134").Value = 8650.5
$ws.Range("M134").Value = -7200.75
$ws.Range("N134").Value = -13720.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 251072.22
$ws.Range("I5").Value = 637
$ws.Range("J5").Value = 304194.84
$ws.Range("K5").Value = 1911
$ws.Range("L5").Value = 912584.52
$ws.Range("M5").Value = -1799
$ws.Range("N5").Value = -912808.52

$ws.Range("H43").Value = 6975
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 6975
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 20925
$ws.Range("N43").Value = -21153

$ws.Range("H92").Value = 801.55554
$ws.Range("I92").Value = 777.6923
$ws.Range("J92").Value = 863.6
$ws.Range("K92").Value = 2333.0769
$ws.Range("L92").Value = 2590.8
$ws.Range("M92").Value = -1085.0769
$ws.Range("N92").Value = -5086.8

$ws.Range("H122").Value = 4248.8335
$ws.Range("I122").Value = 427.64285
$ws.Range("J122").Value = 7592.375
$ws.Range("K122").Value = 3848.78565
$ws.Range("L122").Value = 68331.375
$ws.Range("M122").Value = -1398.78565
$ws.Range("N122").Value = -73231.375

$ws.Range("H132").Value = 1845.5
$ws.Range("I132").Value = 2759.8
$ws.Range("J132").Value = 1493.8462
$ws.Range("K132").Value = 24838.2
$ws.Range("L132").Value = 13444.6158
$ws.Range("M132").Value = -22308.2
$ws.Range("N132").Value = -18504.6158

$ws.Range("H135").Value = 251072.22
$ws.Range("I135").Value = 637
$ws.Range("J135").Value = 304194.84
$ws.Range("K135").Value = 5733
$ws.Range("L135").Value = 2737753.56
$ws.Range("M135").Value = -3198
$ws.Range("N135").Value = -2742823.56

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5724.8237
$ws.Range("I70").Value = 6156.636
$ws.Range("J70").Value = 4933.1665
$ws.Range("K70").Value = 6156.636
$ws.Range("L70").Value = 4933.1665
$ws.Range("M70").Value = -5886.636
$ws.Range("N70").Value = -5473.1665

$ws.Range("H73").Value = 5724.8237
$ws.Range("I73").Value = 6156.636
$ws.Range("J73").Value = 4933.1665
$ws.Range("K73").Value = 6156.636
$ws.Range("L73").Value = 4933.1665
$ws.Range("M73").Value = -5220.636
$ws.Range("N73").Value = -6805.1665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 64497.812
$ws.Range("I7").Value = 68664.336
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 68664.336
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = -68552.336

$ws.Range("H20").Value = 52500
$ws.Range("I20").Value = 5000
$ws.Range("J20").Value = 100000
$ws.Range("K20").Value = 5000
$ws.Range("L20").Value = 100000
$ws.Range("M20").Value = -4774
$ws.Range("N20").Value = -100452

$ws.Range("H40").Value = 26317656
$ws.Range("I40").Value = 34484660
$ws.Range("J40").Value = 1756.1111
$ws.Range("K40").Value = 34484660
$ws.Range("L40").Value = 1756.1111
$ws.Range("M40").Value = -34484524

$ws.Range("H122").Value = 3881506.5
$ws.Range("I122").Value = 5956261
$ws.Range("J122").Value = 1115167.2
$ws.Range("K122").Value = 17868783
$ws.Range("L122").Value = 3345501.6
$ws.Range("M122").Value = -17866333
$ws.Range("N122").Value = -3350401.6

$ws.Range("H126").Value = 64497.812
$ws.Range("I126").Value = 68664.336
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 205993.008
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -203523.008

$ws.Range("H132").Value = 12352128
$ws.Range("I132").Value = 19616610
$ws.Range("J132").Value = 2509.4
$ws.Range("K132").Value = 58849830
$ws.Range("L132").Value = 7528.200000000001
$ws.Range("M132").Value = -58847300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H122").Value = 2208.9
$ws.Range("I122").Value = 1384.875
$ws.Range("J122").Value = 5505
$ws.Range("K122").Value = 4154.625
$ws.Range("L122").Value = 16515
$ws.Range("M122").Value = -1704.625
$ws.Range("N122").Value = -21415

$ws.Range("H126").Value = 1099.4166
$ws.Range("I126").Value = 898.5
$ws.Range("J126").Value = 1501.25
$ws.Range("K126").Value = 2695.5
$ws.Range("L126").Value = 4503.75
$ws.Range("M126").Value = -225.5
$ws.Range("N126").Value = -9443.75
